# Update the cryptocurrency price/volume table on the active sheet.
# Column D ("Price") values that look like a plain decimal number are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (preserving formatting such as trailing zeros, e.g. "4.60") instead of
# silently re-parsing them as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.889.17"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.673.51"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'551.76"
$ws.Range("E5").Value = "  -3.14%  "
$ws.Range("D6").Value = "'158.40"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  -3.73%  "
$ws.Range("E12").Value = "  -7.77%  "
$ws.Range("D13").Value = "3.146.57"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").Value = "'26.15"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "62.777.52"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "'0.0000146"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "2.677.24"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "'11.85"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").Value = "'4.60"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("D20").Value = "'344.41"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "'6.29"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'0.505"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "'63.01"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'8.16"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("E28").Value = "  -5.52%  "
$ws.Range("D29").Value = "'1.36"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'1.93"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").Value = "'167.74"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("D33").Value = "'1.48"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.85"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").Value = "'19.50"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").Value = "'1.78"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").Value = "'348.94"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "'0.960"
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "'6.28"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("D42").Value = "'38.19"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'20.77"
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'20.28"
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("D46").Value = "'0.616"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "'11.02"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'0.0973"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0241"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'128.74"
$ws.Range("E51").Value = "  -4.34%  "
